$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (endesa): Consumo Activa Punta/Llano/Valle. Values look numeric
# ("037", "062", "070") but must be kept as text (leading zeros), so use a
# leading apostrophe to force a text entry, same as typing it in Excel.
$ws.Range("K2").Value = "'037"
$ws.Range("L2").Value = "'062"
$ws.Range("M2").Value = "'070"

# Row 3 (naturgy): Consumo Activa Punta/Llano/Valle
$ws.Range("K3").Value = "29kWh"
$ws.Range("L3").Value = "19kWh"
$ws.Range("M3").Value = "39kWh"

# Row 4 (nexus): Consumo Activa Punta/Llano/Valle
$ws.Range("K4").Value = "6 kWh"
$ws.Range("L4").Value = "7 kWh"
$ws.Range("M4").Value = "4 kWh"
